$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the relative-width / relative-ratio tiers
$ws.Range("S1").Value = "relwidthtiers1"
$ws.Range("T1").Value = "relwidthtiers2"
$ws.Range("U1").Value = "relwidthtiers3"
$ws.Range("V1").Value = "relratiotiers1"
$ws.Range("W1").Value = "relratiotier2"
$ws.Range("X1").Value = "relratiotiers3"

# Add the relative ratio/width formulas for each data row (2 through 15)
for ($r = 2; $r -le 15; $r++) {
    $ws.Range("S$r").Formula = "=J${r}/SUM(J${r}:L${r})"
    $ws.Range("T$r").Formula = "=K${r}/SUM(J${r}:L${r})"
    $ws.Range("U$r").Formula = "=L${r}/(SUM(J${r}:L${r}))"
    $ws.Range("V$r").Formula = "=P${r}/(SUM(P${r}:R${r}))"
    $ws.Range("W$r").Formula = "=Q${r}/SUM(P${r}:R${r})"
    $ws.Range("X$r").Formula = "=R${r}/SUM(P${r}:R${r})"
}

# Update the view: scroll right so column H is the left-most visible
# column (topLeftCell = H1), then restore the last-known selection (T19)
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("T19").Select()
